# Apply the "day one of gwinnett" edit to LookupTable.xlsx:
#  - clear the stray regression-output values that had been pasted into
#    A2/B2 (the lookup table itself only spans columns A:C as headers +
#    the table rows below; A2/B2 weren't meant to hold data here)
#  - move the active selection to F3 to match where the user left off

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Activate()

# Clear the values previously sitting in A2 and B2
$ws.Range("A2:B2").ClearContents()

# Update the active cell selection to F3
$ws.Range("F3").Select()
